# Insert a new slide "@OneToMany (一對多關聯性)" right before the existing
# "JPA的配置" slide (currently slide 16), using the Title+Content layout.
$p = $ppt.ActivePresentation

$s = $p.Slides.Add(16, 2)

# --- Title placeholder ---------------------------------------------------
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "@"
$titleTr.LanguageID = "en-US"

$r = $titleTr.InsertAfter("OneToMany")
$r.LanguageID = "en-US"

$r = $r.InsertAfter(" (")
$r.LanguageID = "en-US"

$r = $r.InsertAfter("一對多關聯性")
$r.LanguageID = "zh-CN"

$r = $r.InsertAfter(")")
$r.LanguageID = "en-US"

# --- Content placeholder --------------------------------------------------
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "單邊的一對多關聯性是指一方有集合屬性，包含多個多方，而多方沒有一方的參考"
$bodyTr.LanguageID = "zh-TW"
